# Apply updated cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.820.37"
$ws.Range("E2").Value = "  +1.37%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.512.06"
$ws.Range("E3").Value = "  +2.50%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.48"
$ws.Range("E5").Value = "  +4.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.98"
$ws.Range("E6").Value = "  +2.79%  "

# Row 7
$ws.Range("E7").Value = "  +0.48%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +2.82%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.512.73"
$ws.Range("E9").Value = "  +1.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  +3.17%  "

# Row 11
$ws.Range("E11").Value = "  -1.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("E12").Value = "  +0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.947.74"
$ws.Range("E14").Value = "  +2.36%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.765.84"
$ws.Range("E15").Value = "  +1.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.37"
$ws.Range("E16").Value = "  +1.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  +2.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.512.94"
$ws.Range("E18").Value = "  +2.38%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  +0.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").Value = "  +2.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.79"
$ws.Range("E21").Value = "  +0.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  +4.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.96"
$ws.Range("E24").Value = "  +4.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.410"
$ws.Range("E25").Value = "  +1.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  +3.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("E29").Value = "  +3.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.29"
$ws.Range("E30").Value = "  +2.36%  "

# Row 31
$ws.Range("E31").Value = "  +4.35%  "

# Row 32
$ws.Range("E32").Value = "  +1.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.32"
$ws.Range("E33").Value = "  +1.16%  "

# Row 34
$ws.Range("E34").Value = "  +0.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +0.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.18"
$ws.Range("E36").Value = "  +1.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("E37").Value = "  -2.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  +0.54%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  +3.11%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.56"
$ws.Range("E40").Value = "  -0.17%  "

# Row 41
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.816"
$ws.Range("E41").Value = "  +6.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.48"
$ws.Range("E42").Value = "  +2.62%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "276.47"
$ws.Range("E43").Value = "  +1.41%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.80"
$ws.Range("E44").Value = "  +10.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.05"
$ws.Range("E45").Value = "  +0.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.593"
$ws.Range("E46").Value = "  +0.75%  "

# Row 47
$ws.Range("E47").Value = "  +2.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0511"
$ws.Range("E48").Value = "  +4.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0219"
$ws.Range("E49").Value = "  +3.84%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.82"
$ws.Range("E50").Value = "  +0.38%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.749.93"
$ws.Range("E51").Value = "  +2.44%  "
